# Remove a set of unused/unnecessary custom paragraph styles from the
# document's style sheet ("Keep unnecessary items out of Word style list").
$d = $word.ActiveDocument

$stylesToRemove = @(
    "internalsectionhead0",
    "BulletedList",
    "BulletedSubsidiaryList",
    "Glossary",
    "GlossaryItem",
    "NumberedList",
    "NumberedSubsidiaryList",
    "UnNumberedList",
    "UnNumberedSubsidiaryList"
)

foreach ($styleName in $stylesToRemove) {
    $style = $d.Styles($styleName)
    if ($style -ne $null) {
        $style.Delete()
    }
}
